# Update countries & provincias Spain
#
# The "Pais" sheet is a COVID-19 dashboard that is re-sorted (descending by
# column B, "Casos totales") every time new daily figures come in. Refreshing
# the numbers below changed the sort order for a handful of countries, which
# is why some rows both get new figures AND a different country label in
# column A (the country that used to occupy that rank moved elsewhere).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4): updated totals -------------------------------
$ws.Cells.Item(4, 2).Value = 1132228   # Casos totales
$ws.Cells.Item(4, 3).Value = 1198      # Nuevos casos
$ws.Cells.Item(4, 5).Value = 904770    # Recuperados
$ws.Cells.Item(4, 7).Value = 39        # Muertes hoy
$ws.Cells.Item(4, 8).Value = 65792     # Muertes

# --- Paises Bajos (row 18) -------------------------------------------------
$ws.Cells.Item(18, 6).Value = 708      # Casos criticos

# --- Oman (row 66) ----------------------------------------------------------
$ws.Cells.Item(66, 5).Value = 1721     # Recuperados
$ws.Cells.Item(66, 7).Value = 1        # Muertes hoy
$ws.Cells.Item(66, 8).Value = 12       # Muertes

# --- Rows 70-73: re-sort + refresh around Ghana/Irak/Uzbekistan/Croacia ----
# Ghana moves up into rank 70 (just after Nigeria) with fresh numbers, and
# Irak, Uzbekistan, Croacia each shift down one rank keeping their own data.

# Row 70: now Ghana (was Irak)
$ws.Cells.Item(70, 1).Value = "Ghana"
$ws.Cells.Item(70, 2).Value = 2169
$ws.Cells.Item(70, 3).Value = 95
$ws.Cells.Item(70, 4).Value = 229
$ws.Cells.Item(70, 5).Value = 1922
$ws.Cells.Item(70, 6).Value = 4
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 18

# Row 71: now Irak (was Uzbekistan)
$ws.Cells.Item(71, 1).Value = "Irak"
$ws.Cells.Item(71, 2).Value = 2153
$ws.Cells.Item(71, 3).Value = 0
$ws.Cells.Item(71, 4).Value = 1414
$ws.Cells.Item(71, 5).Value = 645
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 94

# Row 72: now Uzbekistan (was Croacia)
$ws.Cells.Item(72, 1).Value = "Uzbekistan"
$ws.Cells.Item(72, 2).Value = 2094
$ws.Cells.Item(72, 3).Value = 8
$ws.Cells.Item(72, 4).Value = 1271
$ws.Cells.Item(72, 5).Value = 814
$ws.Cells.Item(72, 6).Value = 8
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 9

# Row 73: now Croacia (was Ghana)
$ws.Cells.Item(73, 1).Value = "Croacia"
$ws.Cells.Item(73, 2).Value = 2088
$ws.Cells.Item(73, 3).Value = 3
$ws.Cells.Item(73, 4).Value = 1463
$ws.Cells.Item(73, 5).Value = 548
$ws.Cells.Item(73, 6).Value = 17
$ws.Cells.Item(73, 7).Value = 2
$ws.Cells.Item(73, 8).Value = 77

# --- Birmania (row 140) ------------------------------------------------------
$ws.Cells.Item(140, 4).Value = 37      # Casos activos
$ws.Cells.Item(140, 5).Value = 108     # Recuperados

# --- Rows 193-194: San Vicente y las Granadinas now ranks just above Namibia
$ws.Cells.Item(193, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(194, 1).Value = "Namibia"

# --- Rows 217-218: San Pedro y Miquelon now ranks just above Comoras -------
$ws.Cells.Item(217, 1).Value = "San Pedro y Miquelon"
$ws.Cells.Item(218, 1).Value = "Comoras"
